$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) New header row (row 5): Promedio (J5), Práctica Efecto Stroop
#    (M5) and Mágico número 7 (N5). K5/L5 keep their existing text.
# -----------------------------------------------------------------
$ws.Range("J5").Value = "Promedio"
$ws.Range("M5").Value = "Práctica Efecto Stroop"
$ws.Range("N5").Value = "Mágico número 7"

# -----------------------------------------------------------------
# 2) Per-row "Práctica Efecto Stroop" (M) / "Mágico número 7" (N)
#    grades, and the new "Promedio" (J) = AVERAGE(K:N) formula.
# -----------------------------------------------------------------
$mValues = @{
   6 = 9.5;  7 = 9.5;  8 = $null; 9 = 7.6; 10 = $null;
  11 = $null; 12 = $null; 13 = $null; 14 = 9.5; 15 = 9;
  16 = 9.1; 17 = 9.5; 18 = $null; 19 = $null; 20 = 9.1;
  21 = 7.6; 22 = 7.6; 23 = 7.6; 24 = $null; 25 = 9.5;
  26 = 9.1; 27 = 9.1; 28 = 9; 29 = 9.5; 30 = 9;
  31 = 9; 32 = 7.6; 33 = 9
}
$nValues = @{
   6 = 9.7;  7 = 9.7;  8 = 8.9; 9 = 9.9; 10 = 8.9;
  11 = 8.9; 12 = $null; 13 = 8.9; 14 = 9.7; 15 = 9.8;
  16 = $null; 17 = 9.7; 18 = 8.9; 19 = 8.9; 20 = $null;
  21 = 9.9; 22 = 9.9; 23 = 9.9; 24 = $null; 25 = 9.7;
  26 = $null; 27 = $null; 28 = 9.8; 29 = 9.7; 30 = 9.8;
  31 = 9.8; 32 = 9.9; 33 = 9.8
}

for ($r = 6; $r -le 33; $r++) {
  $ws.Range("J$r").Formula = "=AVERAGE(K$($r):N$($r))"

  $m = $mValues[$r]
  if ($m -ne $null) {
    $ws.Range("M$r").Value = $m
  }
  $n = $nValues[$r]
  if ($n -ne $null) {
    $ws.Range("N$r").Value = $n
  }
}

# -----------------------------------------------------------------
# 3) Formatting:
#    - J column: plain/default style (nothing special to apply).
#    - L5 header cell should look like K5 (center/center/wrap, no
#      special top-alignment) instead of its old "top aligned" style.
#    - K6:L33 grade columns should wrap their text.
#    - M5:T33 (new Stroop / Mágico número 7 / blank spare columns)
#      should look like K5 (center/center/wrap) too.
# -----------------------------------------------------------------
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("M5:T33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("K6:L33").WrapText = $true

# -----------------------------------------------------------------
# 4) Row heights: every data row (6-33) becomes 46.25pt tall.
# -----------------------------------------------------------------
$ws.Range("A6:A33").RowHeight = 46.25

# -----------------------------------------------------------------
# 5) Refresh the used range / dimension by touching the far corner.
# -----------------------------------------------------------------
$ws.Calculate()
